$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute column G (K) values per updated save_data logic
$ws.Range("G2").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
